# Update the task list on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Insert a new row at row 3 to hold the renamed/relocated target-node task
#    (moves old rows 3..9 down to 4..10)
$ws.Rows.Item(3).Insert()

# Fill the new row 3 with the renamed task (formerly "EnviroDIY Target Node"
# in old row 4), now complete.
$ws.Range("A3").Value = "MyWatershedTarget Node"
$ws.Range("B3").Value = 43647
$ws.Range("C3").Value = "Complete"
$ws.Range("D3").Value = "Develop a target node which posts data to ODM2 Data Portal"

# 2. Remove the old "EnviroDIY Target Node" row, which is now at row 5
#    (old row 4 shifted down by the insert above).
$ws.Rows.Item(5).Delete()

# 3. Mark the "Trend Analysis Node" task (now at row 5) as Active instead of Open.
$ws.Range("C5").Value = "Active"

# 4. Insert a new row after "Trend Analysis Node" (row 5) for the new publisher task.
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "Modular VDAB Publisher"
$ws.Range("B6").Value = 43678
$ws.Range("C6").Value = "Active"
$ws.Range("D6").Value = "Develop a modular publisher library compatible with EnviroDIY modular sensor library."

# 5. Autofit column A now that a wider label was introduced, and update the
#    active selection to reflect where the user left off editing.
$ws.Columns.Item(1).AutoFit()
$ws.Range("A3").Select()
